# redmine # 9229 - Calibration sheet update for GP05MOAS-GL453
# - Moorings!J2: recovery vessel changed from "MV1404" to "CGCS Tully"
# - Moorings!D2: Anchor Launch Date corrected (41868 -> 41880, i.e. 17-Aug-14 -> 29-Aug-14)
# - Asset_Cal_Info!F6: CC_angular_resolution value corrected (1.13 -> 1.096)
# Edited cells are highlighted in blue font (and F6 additionally gets a
# yellow fill) to flag the correction, matching the author's markup.

$wb = $excel.ActiveWorkbook

$moorings = $wb.Worksheets.Item("Moorings")
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet -------------------------------------------------------

# J2: recovering vessel name
$jCell = $moorings.Range("J2")
$jCell.Value = "CGCS Tully"
$jCell.Font.Color = 16711680

# D2: anchor launch date
$dCell = $moorings.Range("D2")
$dCell.Value = 41880
$dCell.Font.Color = 16711680

# --- Asset_Cal_Info sheet --------------------------------------------------

# F6: CC_angular_resolution calibration coefficient value
$fCell = $assetCal.Range("F6")
$fCell.Value = 1.096
$fCell.Font.Color = 16711680
$fCell.Interior.Color = 65535
